# Update master to output generated at c986bee
$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-06 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-10-07 Monday", 2)

# Update the division problems in the table.
# Using direct cell addressing (row, col) avoids any ambiguity from
# duplicate/overlapping text values between old and new content.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "30÷6=5, 0"
$t.Cell(1, 2).Range.Text  = "22÷4=5, 2"
$t.Cell(1, 3).Range.Text  = "54÷4=13, 2"
$t.Cell(1, 4).Range.Text  = "66÷2=33, 0"
$t.Cell(1, 5).Range.Text  = "58÷6=9, 4"

$t.Cell(5, 1).Range.Text  = "92÷8=11, 4"
$t.Cell(5, 2).Range.Text  = "90÷7=12, 6"
$t.Cell(5, 3).Range.Text  = "88÷9=9, 7"
$t.Cell(5, 4).Range.Text  = "10÷4=2, 2"
$t.Cell(5, 5).Range.Text  = "49÷9=5, 4"

$t.Cell(9, 1).Range.Text  = "35÷4=8, 3"
$t.Cell(9, 2).Range.Text  = "41÷2=20, 1"
$t.Cell(9, 3).Range.Text  = "57÷4=14, 1"
$t.Cell(9, 4).Range.Text  = "41÷8=5, 1"
$t.Cell(9, 5).Range.Text  = "63÷4=15, 3"

$t.Cell(13, 1).Range.Text = "85÷2=42, 1"
$t.Cell(13, 2).Range.Text = "72÷8=9, 0"
$t.Cell(13, 3).Range.Text = "97÷5=19, 2"
$t.Cell(13, 4).Range.Text = "74÷7=10, 4"
$t.Cell(13, 5).Range.Text = "93÷2=46, 1"

$t.Cell(17, 1).Range.Text = "75÷9=8, 3"
$t.Cell(17, 2).Range.Text = "12÷4=3, 0"
$t.Cell(17, 3).Range.Text = "49÷4=12, 1"
$t.Cell(17, 4).Range.Text = "64÷8=8, 0"
$t.Cell(17, 5).Range.Text = "38÷2=19, 0"

Write-Host "All replacements applied"
